# Generate Report for Handoff
# Refreshes the generated handoff-report identifiers, content hashes and
# timestamps across the Overview / zh-cn / de-de sheets (as produced by a
# fresh run of the localization handoff report generator).

$wb = $excel.ActiveWorkbook

$oldId   = "a008c084-d8e2-4976-83cf-2fced86014e0"
$newId   = "284d4163-7e09-49d3-b46d-ec46485f8e74"
$oldHash = "5441fffcc9a66c978ec2db0cb255d5af48b72ec0"
$newHash = "3b122179af5fe3bd344b553a9db129a64cc38e5f"

$newMdName    = "$newId.md"
$newZhXlfName = "$newId.$newHash.zh-cn.xlf"
$newDeXlfName = "$newId.$newHash.de-de.xlf"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = "2016-50-15 03:50:13"

foreach ($link in $wsOverview.Hyperlinks) {
    if ($link.Range.Address() -eq '$A$2') {
        $link.TextToDisplay = $newMdName
    }
}

# ---- zh-cn sheet ----
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlfName
$wsZh.Range("E2").Value = "2016-03-15 03:50:03"

foreach ($link in $wsZh.Hyperlinks) {
    if ($link.Range.Address() -eq '$A$2') {
        $link.TextToDisplay = $newMdName
    } elseif ($link.Range.Address() -eq '$D$2') {
        $link.TextToDisplay = $newZhXlfName
    }
}

# ---- de-de sheet ----
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlfName
$wsDe.Range("E2").Value = "2016-03-15 03:50:13"

foreach ($link in $wsDe.Hyperlinks) {
    if ($link.Range.Address() -eq '$A$2') {
        $link.TextToDisplay = $newMdName
    } elseif ($link.Range.Address() -eq '$D$2') {
        $link.TextToDisplay = $newDeXlfName
    }
}
